$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Egf"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07285266666666666
$ws.Range("H2").Value = 0.218558
$ws.Range("I2").Value = 0.05584899373277382
$ws.Range("J2").Value = 0.05584899373277381
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.020000333333333
$ws.Range("N2").Value = 9.060001
$ws.Range("O2").Value = 0.291481777372034
$ws.Range("P2").Value = 0.291481777372034
$ws.Range("Q2").Value = 0.2200150776175555
$ws.Range("R2").Value = 1.980135698558
$ws.Range("S2").Value = 0.0162789639576685
$ws.Range("T2").Value = 0.0162789639576685

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Egf"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07285266666666666
$ws.Range("H3").Value = 0.218558
$ws.Range("I3").Value = 0.05584899373277382
$ws.Range("J3").Value = 0.05584899373277381
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.076388666666666
$ws.Range("N3").Value = 12.229166
$ws.Range("O3").Value = 0.3934413518781783
$ws.Range("P3").Value = 0.3934413518781784
$ws.Range("Q3").Value = 0.2969757847364444
$ws.Range("R3").Value = 2.672782062628
$ws.Range("S3").Value = 0.02197330359525844
$ws.Range("T3").Value = 0.02197330359525844

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Egf"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07285266666666666
$ws.Range("H4").Value = 0.218558
$ws.Range("I4").Value = 0.05584899373277382
$ws.Range("J4").Value = 0.05584899373277381
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.229698
$ws.Range("N4").Value = 9.689094000000001
$ws.Range("O4").Value = 0.311721195201271
$ws.Range("P4").Value = 0.3117211952012711
$ws.Range("Q4").Value = 0.235292111828
$ws.Range("R4").Value = 2.117629006452
$ws.Range("S4").Value = 0.01740931507716855
$ws.Range("T4").Value = 0.01740931507716855

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Egf"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07285266666666666
$ws.Range("H5").Value = 0.218558
$ws.Range("I5").Value = 0.05584899373277382
$ws.Range("J5").Value = 0.05584899373277381
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03476766666666667
$ws.Range("N5").Value = 0.104303
$ws.Range("O5").Value = 0.003355675548516525
$ws.Range("P5").Value = 0.003355675548516525
$ws.Range("Q5").Value = 0.002532917230444444
$ws.Range("R5").Value = 0.022796255074
$ws.Range("S5").Value = 0.0001874111026783217
$ws.Range("T5").Value = 0.0001874111026783217

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Egf"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5726236666666666
$ws.Range("H6").Value = 1.717871
$ws.Range("I6").Value = 0.438974399073536
$ws.Range("J6").Value = 0.438974399073536
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.020000333333333
$ws.Range("N6").Value = 9.060001
$ws.Range("O6").Value = 0.291481777372034
$ws.Range("P6").Value = 0.291481777372034
$ws.Range("Q6").Value = 1.729323664207889
$ws.Range("R6").Value = 15.563912977871
$ws.Range("S6").Value = 0.1279530380627748
$ws.Range("T6").Value = 0.1279530380627748

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Egf"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5726236666666666
$ws.Range("H7").Value = 1.717871
$ws.Range("I7").Value = 0.438974399073536
$ws.Range("J7").Value = 0.438974399073536
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.076388666666666
$ws.Range("N7").Value = 12.229166
$ws.Range("O7").Value = 0.3934413518781783
$ws.Range("P7").Value = 0.3934413518781784
$ws.Range("Q7").Value = 2.334236625065111
$ws.Range("R7").Value = 21.008129625586
$ws.Range("S7").Value = 0.172710681011403
$ws.Range("T7").Value = 0.172710681011403

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Egf"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5726236666666666
$ws.Range("H8").Value = 1.717871
$ws.Range("I8").Value = 0.438974399073536
$ws.Range("J8").Value = 0.438974399073536
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.229698
$ws.Range("N8").Value = 9.689094000000001
$ws.Range("O8").Value = 0.311721195201271
$ws.Range("P8").Value = 0.3117211952012711
$ws.Range("Q8").Value = 1.849401510986
$ws.Range("R8").Value = 16.644613598874
$ws.Range("S8").Value = 0.1368376243419624
$ws.Range("T8").Value = 0.1368376243419624

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Egf"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5726236666666666
$ws.Range("H9").Value = 1.717871
$ws.Range("I9").Value = 0.438974399073536
$ws.Range("J9").Value = 0.438974399073536
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.03476766666666667
$ws.Range("N9").Value = 0.104303
$ws.Range("O9").Value = 0.003355675548516525
$ws.Range("P9").Value = 0.003355675548516525
$ws.Range("Q9").Value = 0.01990878876811111
$ws.Range("R9").Value = 0.179179098913
$ws.Range("S9").Value = 0.0014730556573958
$ws.Range("T9").Value = 0.0014730556573958

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Egf"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.594248
$ws.Range("H10").Value = 1.782744
$ws.Range("I10").Value = 0.4555516544035914
$ws.Range("J10").Value = 0.4555516544035914
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.020000333333333
$ws.Range("N10").Value = 9.060001
$ws.Range("O10").Value = 0.291481777372034
$ws.Range("P10").Value = 0.291481777372034
$ws.Range("Q10").Value = 1.794629158082667
$ws.Range("R10").Value = 16.151662422744
$ws.Range("S10").Value = 0.1327850059103294
$ws.Range("T10").Value = 0.1327850059103294

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Egf"
$ws.Range("C11").Value = "Erbb2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.594248
$ws.Range("H11").Value = 1.782744
$ws.Range("I11").Value = 0.4555516544035914
$ws.Range("J11").Value = 0.4555516544035914
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.076388666666666
$ws.Range("N11").Value = 12.229166
$ws.Range("O11").Value = 0.3934413518781783
$ws.Range("P11").Value = 0.3934413518781784
$ws.Range("Q11").Value = 2.422385812389333
$ws.Range("R11").Value = 21.801472311504
$ws.Range("S11").Value = 0.1792328587588897
$ws.Range("T11").Value = 0.1792328587588897

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Egf"
$ws.Range("C12").Value = "Erbb2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.594248
$ws.Range("H12").Value = 1.782744
$ws.Range("I12").Value = 0.4555516544035914
$ws.Range("J12").Value = 0.4555516544035914
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.229698
$ws.Range("N12").Value = 9.689094000000001
$ws.Range("O12").Value = 0.311721195201271
$ws.Range("P12").Value = 0.3117211952012711
$ws.Range("Q12").Value = 1.919241577104
$ws.Range("R12").Value = 17.273174193936
$ws.Range("S12").Value = 0.1420051061866039
$ws.Range("T12").Value = 0.1420051061866039

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Egf"
$ws.Range("C13").Value = "Erbb2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.594248
$ws.Range("H13").Value = 1.782744
$ws.Range("I13").Value = 0.4555516544035914
$ws.Range("J13").Value = 0.4555516544035914
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03476766666666667
$ws.Range("N13").Value = 0.104303
$ws.Range("O13").Value = 0.003355675548516525
$ws.Range("P13").Value = 0.003355675548516525
$ws.Range("Q13").Value = 0.02066061638133334
$ws.Range("R13").Value = 0.185945547432
$ws.Range("S13").Value = 0.001528683547768382
$ws.Range("T13").Value = 0.001528683547768382

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Egf"
$ws.Range("C14").Value = "Erbb2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.06473366666666668
$ws.Range("H14").Value = 0.194201
$ws.Range("I14").Value = 0.04962495279009878
$ws.Range("J14").Value = 0.04962495279009878
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.020000333333333
$ws.Range("N14").Value = 9.060001
$ws.Range("O14").Value = 0.291481777372034
$ws.Range("P14").Value = 0.291481777372034
$ws.Range("Q14").Value = 0.1954956949112223
$ws.Range("R14").Value = 1.759461254201
$ws.Range("S14").Value = 0.01446476944126127
$ws.Range("T14").Value = 0.01446476944126127

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Egf"
$ws.Range("C15").Value = "Erbb2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.06473366666666668
$ws.Range("H15").Value = 0.194201
$ws.Range("I15").Value = 0.04962495279009878
$ws.Range("J15").Value = 0.04962495279009878
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.076388666666666
$ws.Range("N15").Value = 12.229166
$ws.Range("O15").Value = 0.3934413518781783
$ws.Range("P15").Value = 0.3934413518781784
$ws.Range("Q15").Value = 0.2638795851517778
$ws.Range("R15").Value = 2.374916266366
$ws.Range("S15").Value = 0.01952450851262724
$ws.Range("T15").Value = 0.01952450851262724

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Egf"
$ws.Range("C16").Value = "Erbb2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.06473366666666668
$ws.Range("H16").Value = 0.194201
$ws.Range("I16").Value = 0.04962495279009878
$ws.Range("J16").Value = 0.04962495279009878
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.229698
$ws.Range("N16").Value = 9.689094000000001
$ws.Range("O16").Value = 0.311721195201271
$ws.Range("P16").Value = 0.3117211952012711
$ws.Range("Q16").Value = 0.2090701937660001
$ws.Range("R16").Value = 1.881631743894
$ws.Range("S16").Value = 0.01546914959553624
$ws.Range("T16").Value = 0.01546914959553624

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Egf"
$ws.Range("C17").Value = "Erbb2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.06473366666666668
$ws.Range("H17").Value = 0.194201
$ws.Range("I17").Value = 0.04962495279009878
$ws.Range("J17").Value = 0.04962495279009878
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.03476766666666667
$ws.Range("N17").Value = 0.104303
$ws.Range("O17").Value = 0.003355675548516525
$ws.Range("P17").Value = 0.003355675548516525
$ws.Range("Q17").Value = 0.002250638544777778
$ws.Range("R17").Value = 0.020255746903
$ws.Range("S17").Value = 0.0001665252406740214
$ws.Range("T17").Value = 0.0001665252406740214

Write-Host "Applied new TPM data: 16 rows x 20 cols"